$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation (2022-03-21) needs to be inserted at row 163,
# pushing the existing rows 163-194 down to 164-195.
$ws.Rows(163).Insert()

$ws.Cells.Item(163, 1).Value = 4
$ws.Cells.Item(163, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(163, 3).Value = "Los Lagos"
$ws.Cells.Item(163, 4).Value = 44641
$ws.Cells.Item(163, 5).Value = 10
$ws.Cells.Item(163, 6).Value = 100112039
$ws.Cells.Item(163, 7).Value = "Ciboulette"
$ws.Cells.Item(163, 8).Value = "Sin especificar"
$ws.Cells.Item(163, 9).Value = "Primera"
$ws.Cells.Item(163, 10).Value = 80
$ws.Cells.Item(163, 11).Value = 6000
$ws.Cells.Item(163, 12).Value = 6000
$ws.Cells.Item(163, 13).Value = 6000
$ws.Cells.Item(163, 14).Value = "`$/docena de atados"
$ws.Cells.Item(163, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(163, 16).Value = 2000
$ws.Cells.Item(163, 17).Value = 3
$ws.Cells.Item(163, 18).Value = "Hortaliza"
